# "Files must now only contain account, debit and credit column"
#
# The Trial Balance sheet previously carried extra descriive columns that
# aren't part of the account/debit/credit shape the report should expose:
#   - E2:E8  the title block (company name, blanks, "Trial Balance", period)
#   - H11:H54 / H56  a "category" label column (and its balance-check formula)
#
# Clear their contents (keeping the existing cell formatting/styles intact)
# so the sheet only shows Account / YTD Debit / YTD Credit going forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title block above the table (company name / blank rows / report title / period)
$ws.Range("E2:E8").ClearContents()

# Drop the extra category column next to the totals row, including the
# G56-F56 balance-check formula in H56.
$ws.Range("H11:H56").ClearContents()

# Reflect the author's new focus: the cursor/selection moves onto the
# now-unused column H, and the view is scrolled down towards the bottom
# of the table.
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H1:H1048576").Select()
